$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.43407928943634
$ws.Range("B1").Value = 1.634546637535095
$ws.Range("C1").Value = 2.033814191818237
$ws.Range("D1").Value = 2.103904247283936
$ws.Range("E1").Value = 1.54181444644928
